# Prototype BMB standards registry
# Adds 96 additional package rows ("pack1".."pack96") to the "packages"
# worksheet of the EMX all-datatypes test fixture.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("packages")

for ($i = 1; $i -le 96; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 1).Value = "pack$i"
}

$ws.Activate()
$ws.Range("A2").Select()
